$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

$ws.Range('D2').Value = '60.600.41'
$ws.Range('E2').Value = '  +4.74%  '
$ws.Range('D3').Value = '2.353.81'
$ws.Range('E3').Value = '  +3.02%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '547.55'
$ws.Range('E5').Value = '  +2.78%  '
Set-TextValue 'D6' '132.61'
$ws.Range('E6').Value = '  +1.77%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +1.59%  '
$ws.Range('D9').Value = '2.351.66'
$ws.Range('E9').Value = '  +2.80%  '
$ws.Range('E10').Value = '  +2.16%  '
Set-TextValue 'D11' '5.51'
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('E12').Value = '  +1.12%  '
Set-TextValue 'D13' '0.335'
$ws.Range('E13').Value = '  +1.86%  '
Set-TextValue 'D14' '24.01'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('D15').Value = '2.772.18'
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').Value = '60.540.35'
$ws.Range('E16').Value = '  +4.72%  '
Set-TextValue 'D17' '0.0000134'
$ws.Range('E17').Value = '  +2.24%  '
$ws.Range('D18').Value = '2.347.99'
$ws.Range('E18').Value = '  +2.37%  '
Set-TextValue 'D19' '10.71'
$ws.Range('E19').Value = '  +2.37%  '
$ws.Range('E20').Value = '  -0.67%  '
Set-TextValue 'D21' '6.88'
$ws.Range('E21').Value = '  +8.35%  '
Set-TextValue 'D22' '315.15'
$ws.Range('E22').Value = '  +1.15%  '
Set-TextValue 'D23' '1.00'
$ws.Range('E23').Value = '  +0.06%  '
Set-TextValue 'D24' '63.40'
$ws.Range('E24').Value = '  +1.75%  '
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('E26').Value = '  +0.08%  '
Set-TextValue 'D27' '7.95'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  +5.93%  '
$ws.Range('E29').Value = '  +3.31%  '
Set-TextValue 'D30' '171.98'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0734'
$ws.Range('E31').Value = '  +2.82%  '
$ws.Range('B32').Value = 'SuiNetwork'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue 'D32' '1.15'
$ws.Range('E32').Value = '  +10.59%  '
Set-TextValue 'D33' '5.93'
$ws.Range('E33').Value = '  +3.51%  '
Set-TextValue 'D34' '1.42'
$ws.Range('E34').Value = '  +15.60%  '
Set-TextValue 'D35' '0.382'
$ws.Range('E35').Value = '  +1.12%  '
Set-TextValue 'D36' '18.07'
$ws.Range('E36').Value = '  +2.18%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  -0.03%  '
Set-TextValue 'D39' '4.17'
$ws.Range('E39').Value = '  +7.63%  '
Set-TextValue 'D40' '316.25'
$ws.Range('E40').Value = '  +10.61%  '
Set-TextValue 'D41' '38.22'
$ws.Range('E41').Value = '  +0.24%  '
Set-TextValue 'D42' '1.54'
$ws.Range('E42').Value = '  +3.96%  '
Set-TextValue 'D43' '142.65'
$ws.Range('E43').Value = '  +1.09%  '
$ws.Range('E44').Value = '  +2.26%  '
$ws.Range('E45').Value = '  +1.25%  '
Set-TextValue 'D46' '19.30'
$ws.Range('E46').Value = '  +7.11%  '
Set-TextValue 'D47' '0.0499'
$ws.Range('E47').Value = '  +1.22%  '
Set-TextValue 'D48' '0.562'
$ws.Range('E48').Value = '  +1.55%  '
Set-TextValue 'D49' '0.0215'
$ws.Range('E49').Value = '  +2.61%  '
Set-TextValue 'D50' '11.04'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').Value = '0.0₆0209'
$ws.Range('E51').Value = '  +4.00%  '
